# StorageType.xlsx - "added date and automation status"
# Adds an "Automated" (yes) + "DateOfAutomation" (3/4/2011) value pair to every
# data row on Sheet1 (rows 2-18), sizes the two new columns, and nudges the
# selection over to the first newly-populated cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- size the new C/D columns (existing A/B columns keep the sheet default) ---
$ws.Columns.Item(3).ColumnWidth = 13.9
$ws.Columns.Item(4).ColumnWidth = 13.2

# --- date value used for every row: 3/4/2011 (serial 40606) ---
$automationDate = 40606

# --- stamp every data row (2-18) with Automated = "yes" / DateOfAutomation ---
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = "yes"
    $ws.Cells.Item($r, 4).Value = $automationDate
    $ws.Cells.Item($r, 4).NumberFormat = "mm-dd-yy"
}

# Rows 2-5,7-12,14-18 now fit on a single (default) line again; let Excel
# drop their stale explicit height. Rows 6 & 13 hold the longest label and
# still wrap onto two lines.
for ($r = 2; $r -le 18; $r++) {
    $ws.Rows.Item($r).AutoFit()
}
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30

# --- move the selection to the first newly-entered cell ---
$ws.Range("C2").Select()
